$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "基金受益憑證" (fund) sheet: row 1 currently holds a duplicate of the
#    data row instead of column headers, and the data row (row 2) is missing
#    the standard trailing metadata columns that every other sheet has
#    (property_category, category, date, legislator_name, legislator_id,
#    source_file, index). Extend the used range from B1:H2 to B1:O2, turn
#    row 1 into proper headers and fill in the missing metadata on row 2.
# ---------------------------------------------------------------------------
$wsFund = $wb.Worksheets.Item("基金受益憑證")

# Grow the formatted range so the new columns inherit the existing look
# (bold/bordered header style on row 1, plain style on row 2).
$wsFund.Range("B1:H2").Copy()
$wsFund.Range("I1:O2").PasteSpecial(-4122)

# Row 1 becomes the header row.
$wsFund.Range("B1").Value = "name"
$wsFund.Range("C1").Value = "owner"
$wsFund.Range("D1").Value = "dealer"
$wsFund.Range("E1").Value = "quantity"
$wsFund.Range("F1").Value = "face_value"
$wsFund.Range("G1").Value = "currency"
$wsFund.Range("H1").Value = "total"
$wsFund.Range("I1").Value = "property_category"
$wsFund.Range("J1").Value = "category"
$wsFund.Range("K1").Value = "date"
$wsFund.Range("L1").Value = "legislator_name"
$wsFund.Range("M1").Value = "legislator_id"
$wsFund.Range("N1").Value = "source_file"
$wsFund.Range("O1").Value = "index"

# Row 2 keeps its existing name/owner/dealer/quantity/face_value/currency/total
# values and just gains the metadata columns.
$wsFund.Range("I2").Value = "fund"
$wsFund.Range("J2").Value = "normal"

# "2011-11-21" looks like a date, so Excel would silently coerce a plain
# assignment into a date serial; force it to stay text first, then restore
# the shared header/data look from the already-populated sibling column.
$wsFund.Range("K2").NumberFormat = "@"
$wsFund.Range("K2").Value = "2011-11-21"
$wsFund.Range("C2").Copy()
$wsFund.Range("K2").PasteSpecial(-4122)

$wsFund.Range("L2").Value = "羅淑蕾"
$wsFund.Range("M2").Value = 1638
$wsFund.Range("N2").Value = "tmpa0031"
$wsFund.Range("O2").Value = 123

# ---------------------------------------------------------------------------
# 2) "其他有價證券" sheet never held real data - rows 1-3 are mangled
#    leftovers of the "(9) antiques / other valuable property" section
#    header from the source PDF, not an actual holding. Drop the sheet.
# ---------------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("其他有價證券")
$wsOther.Delete()

